# Update "想去人数" (F column) counts for a few exhibition events.
# Affects both the "展览" sheet (sheet1) and the "全部类型" sheet (sheet4),
# since the latter aggregates rows from all the category sheets.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAll        = $wb.Worksheets.Item("全部类型")

# 展览 sheet: rows 2, 4, 5, 11 -> column F
$wsExhibition.Range("F2").Value  = 303
$wsExhibition.Range("F4").Value  = 8061
$wsExhibition.Range("F5").Value  = 5877
$wsExhibition.Range("F11").Value = 405

# 全部类型 sheet: rows 2, 4, 5, 14 -> column F
$wsAll.Range("F2").Value  = 303
$wsAll.Range("F4").Value  = 8061
$wsAll.Range("F5").Value  = 5877
$wsAll.Range("F14").Value = 405
